$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @(901, 16, 15, 45, 60, 60)
    3 = @(301, 6, 45, 30, 60, 45)
    4 = @(801, 3, 67, 65, 52, 45)
    5 = @(401, 9, 48, 67, 75, 45)
    6 = @(1201, 2, 10, 10, 10, 10)
    7 = @(1202, 2, 10, 10, 10, 10)
    8 = @(201, 9, 30, 15, 45, 30)
    9 = @(1203, 3, 15, 15, 15, 15)
    10 = @(902, 1, 0, 0, 0, 0)
    11 = @(501, 9, 52, 30, 75, 45)
    12 = @(701, 3, 90, 45, 97, 15)
    13 = @(101, 9, 30, 15, 60, 15)
    14 = @(1001, 18, 30, 75, 60, 72)
    16 = @(502, 0, 4, 0, 0, 0)
    17 = @(802, 0, 4, 5, 4, 0)
    18 = @(1, 0, 2, 2, 2, 2)
    19 = @(2, 0, 2, 2, 2, 2)
    22 = @(602, 0, 0, 4, 0, 9)
    23 = @(402, 0, 0, 4, 0, 0)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
